# Update "想去人数" (number of people interested) figures that changed
# between the two scrapes reflected in this commit.
#
# Sheet "展览"   (Exhibitions): row 2 -> F2: 416 -> 417 ; row 3 -> F3: 2538 -> 2561
# Sheet "全部类型" (All types) : row 2 -> F2: 416 -> 417 ; row 7 -> F7: 2538 -> 2561

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 417
$wsExhibit.Range("F3").Value = 2561

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 417
$wsAll.Range("F7").Value = 2561
